# "Generate Report for Archive"
#
# The localization-status report is regenerated: the status of the
# 8334310f-... source file moves from "Ready for handoff" to
# "In Translation" on every sheet that lists it (Overview, zh-cn, de-de),
# and the (now shorter) status column is re-sized to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns are E and F.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Per-locale sheets: Status column is C.
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Shrink the status columns to match the shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
